$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new/updated rows: stock ticker and price (price kept as text,
# matching the original sheet's inline-string representation of "$x.xx").
$data = @(
    @(0,  "CHR:TSE",      '$3.01'),
    @(1,  "HUT:TSE",      '$30.85'),
    @(2,  "GWO:TSE",      '$47.50'),
    @(3,  "CGX:TSE",      '$12.16'),
    @(4,  "BB:TSE",       '$5.51'),
    @(5,  "VRN:TSE",      '$7.17'),
    @(6,  "TLRY:TSE",     '$1.96'),
    @(7,  "PXT:TSE",      '$13.67'),
    @(8,  "AC:TSE",       '$22.25'),
    @(9,  "SU:TSE",       '$50.71'),
    @(10, "BEN:NYSE",     '$20.11'),
    @(11, "ACB:TSE",      '$6.03'),
    @(12, "OGI:TSE",      '$2.27'),
    @(13, "POU:TSE",      '$31.22'),
    @(14, "ASM:TSE",      '$1.25'),
    @(15, "QCLN:NASDAQ",  '$34.47'),
    @(16, "AMAT:NASDAQ",  '$163.64'),
    @(17, "OGI:TSE",      '$2.27'),
    @(18, "POU:TSE",      '$31.22'),
    @(19, "QQC:TSE",      '$36.15')
)

# Extend the styled (bold / bordered / centered) look of A2 down through A21
# before filling in the numbers, so every row-index cell in column A picks
# up the same style used by the existing A2 cell.
$ws.Range("A2").Copy()
$ws.Range("A3:A21").PasteSpecial(-4122)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]

    $ws.Cells.Item($row, 2).Value = $item[1]

    # Prices look like currency ("$3.01"), which Excel would otherwise
    # auto-convert to a numeric currency value. Force them to be stored
    # as plain text, then strip the resulting formatting override so the
    # cell keeps the sheet's default (unstyled) look, matching the rest
    # of the table.
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 3).Style = "Normal"

    $row++
}
